$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new "Label" header in H1, copying the header style from G1
$ws.Range("H1").Value = "Label"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# Updated D/E values for the 100-iteration block (rows 2-11) and H column
# labels (0 = Control, 1 = MDD) for every data row.

$ws.Range("D2").Value = 0.3702967762937955
$ws.Range("E2").Value = 0.3702967762937955
$ws.Range("H2").Value = 0

$ws.Range("H3").Value = 0

$ws.Range("D4").Value = 0.4017428440661174
$ws.Range("E4").Value = 0.4017428440661174
$ws.Range("H4").Value = 0

$ws.Range("D5").Value = 0.2875422608624412
$ws.Range("E5").Value = 0.2875422608624412
$ws.Range("H5").Value = 0

$ws.Range("D6").Value = 0.684537825908332
$ws.Range("E6").Value = 0.684537825908332
$ws.Range("H6").Value = 0

$ws.Range("D7").Value = 0.3733273479538537
$ws.Range("E7").Value = 0.6266726520461463
$ws.Range("H7").Value = 1

$ws.Range("D8").Value = 0.5740043159402053
$ws.Range("E8").Value = 0.4259956840597947
$ws.Range("H8").Value = 1

$ws.Range("D9").Value = 0.5202069637837331
$ws.Range("E9").Value = 0.4797930362162669
$ws.Range("H9").Value = 1

$ws.Range("D10").Value = 0.3797952983583848
$ws.Range("E10").Value = 0.6202047016416152
$ws.Range("H10").Value = 1

$ws.Range("D11").Value = 0.2935776902317846
$ws.Range("E11").Value = 0.7064223097682154
$ws.Range("F11").Value = 0.7069457769393921
$ws.Range("H11").Value = 1

# 200-iteration block (rows 12-21): D/E/F/G unchanged, only the new H
# (Label) column is populated.

$ws.Range("H12").Value = 0
$ws.Range("H13").Value = 0
$ws.Range("H14").Value = 0
$ws.Range("H15").Value = 0
$ws.Range("H16").Value = 0
$ws.Range("H17").Value = 1
$ws.Range("H18").Value = 1
$ws.Range("H19").Value = 1
$ws.Range("H20").Value = 1
$ws.Range("H21").Value = 1
